# Update cryptos list cell values.
# Each text cell is set using a leading apostrophe to force Excel to
# treat the value as literal text (preventing numeric auto-conversion
# for values such as "0.626" or "68.158.56"), then the cell style is
# reset back to "Normal" so no stray quotePrefix formatting remains.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.158.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.90%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.915.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'484.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.71%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'146.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.69%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.53%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.728"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.71%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.02%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0000355"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'42.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.67%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'10.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.15%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.528.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.02%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'14.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.916.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.40%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.18%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'19.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -2.34%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'68.263.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.68%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'447.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.39%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'14.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.31%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.60%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'88.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'11.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +14.36%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +13.42%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.26%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'38.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.51%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +3.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Cosmos"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'13.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.98%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Bittensor"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'689.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -6.68%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.95%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.87%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0₃0919"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +20.49%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'41.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'59.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.86%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +5.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.150"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +18.53%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0478"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.368"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +9.31%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.94%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.06%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.84%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.73%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'146.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.02%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.59%  "
$ws.Range("E51").Style = "Normal"
